$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3466
$ws.Range("I19").Value = 2554.6667
$ws.Range("J19").Value = 4559.6
$ws.Range("K19").Value = 2554.6667
$ws.Range("L19").Value = 4559.6
$ws.Range("M19").Value = -2379.6667
$ws.Range("N19").Value = -4909.6

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H52").Value = 264.22858
$ws.Range("I52").Value = 125.28571
$ws.Range("K52").Value = 375.85713
$ws.Range("M52").Value = -215.85713

$ws.Range("H92").Value = 30303656
$ws.Range("I92").Value = 465.5
$ws.Range("K92").Value = 465.5
$ws.Range("M92").Value = 782.5

$ws.Range("H106").Value = 3182.1428
$ws.Range("I106").Value = 3266.25
$ws.Range("K106").Value = 3266.25
$ws.Range("M106").Value = -2635.25

$ws.Range("H132").Value = 1850.9565
$ws.Range("I132").Value = 1860.15
$ws.Range("K132").Value = 5580.450000000001
$ws.Range("M132").Value = -3050.450000000001

$ws.Range("H138").Value = 1045339.9
$ws.Range("I138").Value = 2482.1035
$ws.Range("K138").Value = 7446.310500000001
$ws.Range("M138").Value = -2306.310500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20002476
$ws.Range("I2").Value = 1169.625
$ws.Range("J2").Value = 100007704
$ws.Range("K2").Value = 1169.625
$ws.Range("L2").Value = 100007704
$ws.Range("M2").Value = -1056.625
$ws.Range("N2").Value = -100007930

$ws.Range("H26").Value = 3625
$ws.Range("I26").Value = 3625
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3625
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -3295
$ws.Range("N26").ClearContents()

$ws.Range("H32").Value = 2480763.5
$ws.Range("I32").Value = 2990995.5
$ws.Range("J32").Value = 38938.43
$ws.Range("K32").Value = 2990995.5
$ws.Range("L32").Value = 38938.43
$ws.Range("M32").Value = -2990708.5
$ws.Range("N32").Value = -39512.43

$ws.Range("H45").Value = 6478.5
$ws.Range("I45").Value = 2357.6
$ws.Range("K45").Value = 2357.6
$ws.Range("M45").Value = -1980.6

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H61").Value = 15153816
$ws.Range("I61").Value = 1478.804
$ws.Range("K61").Value = 1478.804
$ws.Range("M61").Value = -1266.804

$ws.Range("H63").Value = 6078.5835
$ws.Range("I63").Value = 6844.3
$ws.Range("K63").Value = 6844.3
$ws.Range("M63").Value = -6158.3

$ws.Range("H66").Value = 6078.5835
$ws.Range("I66").Value = 6844.3
$ws.Range("K66").Value = 34221.5
$ws.Range("M66").Value = -30789.5

$ws.Range("H74").Value = 34122.594
$ws.Range("I74").Value = 51746.9
$ws.Range("J74").Value = 4748.75
$ws.Range("K74").Value = 51746.9
$ws.Range("L74").Value = 4748.75
$ws.Range("M74").Value = -50872.9
$ws.Range("N74").Value = -6496.75

$ws.Range("H77").Value = 34122.594
$ws.Range("I77").Value = 51746.9
$ws.Range("J77").Value = 4748.75
$ws.Range("K77").Value = 258734.5
$ws.Range("L77").Value = 23743.75
$ws.Range("M77").Value = -254366.5
$ws.Range("N77").Value = -32479.75

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H102").Value = 3934.125
$ws.Range("I102").Value = 3346
$ws.Range("K102").Value = 3346
$ws.Range("M102").Value = -1724

$ws.Range("H110").Value = 12346835
$ws.Range("I110").Value = 1055.9524
$ws.Range("K110").Value = 1055.9524
$ws.Range("M110").Value = 989.0476000000001

$ws.Range("H112").Value = 59343
$ws.Range("J112").Value = 59343
$ws.Range("L112").Value = 59343
$ws.Range("N112").Value = -62297

$ws.Range("H116").Value = 20002476
$ws.Range("I116").Value = 1169.625
$ws.Range("J116").Value = 100007704
$ws.Range("K116").Value = 1169.625
$ws.Range("L116").Value = 100007704
$ws.Range("M116").Value = 1124.375
$ws.Range("N116").Value = -100012292

$ws.Range("H122").Value = 3751.4443
$ws.Range("I122").Value = 2537.5715
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 7612.7145
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5162.7145
$ws.Range("N122").Value = -28900

$ws.Range("H132").Value = 2212.3125
$ws.Range("I132").Value = 995.75806
$ws.Range("K132").Value = 2987.27418
$ws.Range("M132").Value = -457.2741799999999

$ws.Range("H136").Value = 15153816
$ws.Range("I136").Value = 1478.804
$ws.Range("K136").Value = 4436.412
$ws.Range("M136").Value = -1886.412

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20002476
$ws.Range("I3").Value = 1169.625
$ws.Range("J3").Value = 100007704
$ws.Range("K3").Value = 1169.625
$ws.Range("L3").Value = 100007704
$ws.Range("M3").Value = -1055.625
$ws.Range("N3").Value = -100007932

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H64").Value = 25641566
$ws.Range("J64").Value = 627.2
$ws.Range("L64").Value = 627.2
$ws.Range("N64").Value = -1077.2

$ws.Range("H67").Value = 25641566
$ws.Range("J67").Value = 627.2
$ws.Range("L67").Value = 627.2
$ws.Range("N67").Value = -2187.2

$ws.Range("H94").Value = 2530.9312
$ws.Range("I94").Value = 1362.5
$ws.Range("K94").Value = 1362.5
$ws.Range("M94").Value = -911.5

$ws.Range("H105").Value = 3039.1667
$ws.Range("I105").Value = 2553.1428
$ws.Range("K105").Value = 2553.1428
$ws.Range("M105").Value = -806.1428000000001

$ws.Range("H107").Value = 70376530
$ws.Range("I107").Value = 75068184
$ws.Range("J107").Value = 1689
$ws.Range("K107").Value = 75068184
$ws.Range("L107").Value = 1689
$ws.Range("M107").Value = -75066264
$ws.Range("N107").Value = -5529

$ws.Range("H134").Value = 6253992.5
$ws.Range("I134").Value = 11906081
$ws.Range("K134").Value = 35718243
$ws.Range("M134").Value = -35715708

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6920.029
$ws.Range("I31").Value = 3520
$ws.Range("J31").Value = 7486.7
$ws.Range("K31").Value = 3520
$ws.Range("L31").Value = 7486.7
$ws.Range("M31").Value = -3225
$ws.Range("N31").Value = -8076.7

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H34").Value = 6920.029
$ws.Range("I34").Value = 3520
$ws.Range("J34").Value = 7486.7
$ws.Range("K34").Value = 3520
$ws.Range("L34").Value = 7486.7
$ws.Range("M34").Value = -3318
$ws.Range("N34").Value = -7890.7

$ws.Range("H36").Value = 39907.4
$ws.Range("I36").Value = 48
$ws.Range("J36").Value = 49872.25
$ws.Range("K36").Value = 48
$ws.Range("L36").Value = 49872.25
$ws.Range("M36").Value = 340
$ws.Range("N36").Value = -50648.25

$ws.Range("H40").Value = 39907.4
$ws.Range("I40").Value = 48
$ws.Range("J40").Value = 49872.25
$ws.Range("K40").Value = 48
$ws.Range("L40").Value = 49872.25
$ws.Range("M40").Value = 112
$ws.Range("N40").Value = -50192.25

$ws.Range("H41").Value = 40021.668
$ws.Range("J41").Value = 57532.5
$ws.Range("L41").Value = 57532.5
$ws.Range("N41").Value = -58388.5

$ws.Range("H58").Value = 4848.7334
$ws.Range("I58").Value = 3452.24
$ws.Range("K58").Value = 3452.24
$ws.Range("M58").Value = -3249.24

$ws.Range("H94").Value = 2087.9092
$ws.Range("I94").Value = 2273.8333
$ws.Range("J94").Value = 1864.8
$ws.Range("K94").Value = 2273.8333
$ws.Range("L94").Value = 1864.8
$ws.Range("M94").Value = -1822.8333
$ws.Range("N94").Value = -2766.8

$ws.Range("H132").Value = 4733.7837
$ws.Range("I132").Value = 3347.9524
$ws.Range("K132").Value = 10043.8572
$ws.Range("M132").Value = -7513.8572

$ws.Range("H134").Value = 3603.3865
$ws.Range("I134").Value = 2389.3547
$ws.Range("K134").Value = 7168.0641
$ws.Range("M134").Value = -4633.0641

$ws.Range("H136").Value = 4848.7334
$ws.Range("I136").Value = 3452.24
$ws.Range("K136").Value = 10356.72
$ws.Range("M136").Value = -7806.719999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43168180
$ws.Range("J4").Value = 3935661.5
$ws.Range("L4").Value = 11806984.5
$ws.Range("N4").Value = -11807208.5

$ws.Range("H18").Value = 619.0909
$ws.Range("I18").Value = 440.14285
$ws.Range("K18").Value = 1320.42855
$ws.Range("M18").Value = -1151.42855

$ws.Range("H23").Value = 503.72726
$ws.Range("J23").Value = 409.6
$ws.Range("L23").Value = 1228.8
$ws.Range("N23").Value = -1698.8

$ws.Range("H39").Value = 13474.25
$ws.Range("J39").Value = 13474.25
$ws.Range("L39").Value = 40422.75
$ws.Range("N39").Value = -41010.75

$ws.Range("H55").Value = 9106346
$ws.Range("J55").Value = 11128777
$ws.Range("L55").Value = 33386331
$ws.Range("N55").Value = -33386685

$ws.Range("H121").Value = 1924357.8
$ws.Range("I121").Value = 1280.5385
$ws.Range("J121").Value = 3847435
$ws.Range("K121").Value = 3841.6155
$ws.Range("L121").Value = 11542305
$ws.Range("M121").Value = -2531.6155
$ws.Range("N121").Value = -11544925

$ws.Range("H122").Value = 1415165.6
$ws.Range("J122").Value = 543.5
$ws.Range("L122").Value = 4891.5
$ws.Range("N122").Value = -9791.5

$ws.Range("H132").Value = 3839.2888
$ws.Range("I132").Value = 2496.45
$ws.Range("J132").Value = 4913.56
$ws.Range("K132").Value = 22468.05
$ws.Range("L132").Value = 44222.04
$ws.Range("M132").Value = -19938.05
$ws.Range("N132").Value = -49282.04

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7445.7
$ws.Range("I70").Value = 6579
$ws.Range("J70").Value = 10912.5
$ws.Range("K70").Value = 6579
$ws.Range("L70").Value = 10912.5
$ws.Range("M70").Value = -6309
$ws.Range("N70").Value = -11452.5

$ws.Range("H73").Value = 7445.7
$ws.Range("I73").Value = 6579
$ws.Range("J73").Value = 10912.5
$ws.Range("K73").Value = 6579
$ws.Range("L73").Value = 10912.5
$ws.Range("M73").Value = -5643
$ws.Range("N73").Value = -12784.5

$ws.Range("H97").Value = 2311.2856
$ws.Range("I97").Value = 1092.3334
$ws.Range("K97").Value = 1092.3334
$ws.Range("M97").Value = -596.3334

$ws.Range("H113").Value = 2507.2273
$ws.Range("I113").Value = 2107.0833
$ws.Range("J113").Value = 2987.4
$ws.Range("K113").Value = 2107.0833
$ws.Range("L113").Value = 2987.4
$ws.Range("M113").Value = 62.91670000000022
$ws.Range("N113").Value = -7327.4

$ws.Range("H122").Value = 1692518
$ws.Range("J122").Value = 3975.75
$ws.Range("L122").Value = 11927.25
$ws.Range("N122").Value = -16827.25

$ws.Range("H132").Value = 1385.2894
$ws.Range("I132").Value = 711.871
$ws.Range("J132").Value = 4367.5713
$ws.Range("K132").Value = 2135.613
$ws.Range("L132").Value = 13102.7139
$ws.Range("M132").Value = 394.3870000000002
$ws.Range("N132").Value = -18162.7139

$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 180000
$ws.Range("N134").Value = -185070

$ws.Range("H136").Value = 44024.39
$ws.Range("J136").Value = 46976.895
$ws.Range("L136").Value = 140930.685
$ws.Range("N136").Value = -146030.685

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4782.1377
$ws.Range("I7").Value = 4247.05
$ws.Range("K7").Value = 4247.05
$ws.Range("M7").Value = -4135.05

$ws.Range("H16").Value = 995.7
$ws.Range("I16").Value = 884.1667
$ws.Range("J16").Value = 1999.5
$ws.Range("K16").Value = 884.1667
$ws.Range("L16").Value = 1999.5
$ws.Range("M16").Value = -714.1667
$ws.Range("N16").Value = -2339.5

$ws.Range("H22").Value = 1267.6666
$ws.Range("I22").Value = 482.1
$ws.Range("K22").Value = 482.1
$ws.Range("M22").Value = -187.1

$ws.Range("H27").Value = 1267.6666
$ws.Range("I27").Value = 482.1
$ws.Range("K27").Value = 482.1
$ws.Range("M27").Value = -375.1

$ws.Range("H63").Value = 41692.332
$ws.Range("I63").Value = 40077
$ws.Range("J63").Value = 42500
$ws.Range("K63").Value = 40077
$ws.Range("L63").Value = 42500
$ws.Range("M63").Value = -39328
$ws.Range("N63").Value = -43998

$ws.Range("H66").Value = 41692.332
$ws.Range("I66").Value = 40077
$ws.Range("J66").Value = 42500
$ws.Range("K66").Value = 120231
$ws.Range("L66").Value = 127500
$ws.Range("M66").Value = -116487
$ws.Range("N66").Value = -134988

$ws.Range("H100").Value = 3273.3572
$ws.Range("I100").Value = 1442.7
$ws.Range("J100").Value = 7850
$ws.Range("K100").Value = 1442.7
$ws.Range("L100").Value = 7850
$ws.Range("M100").Value = -901.7
$ws.Range("N100").Value = -8932

$ws.Range("H110").Value = 1000000000
$ws.Range("J110").Value = 1000000000
$ws.Range("L110").Value = 1000000000
$ws.Range("N110").Value = -1000008180

$ws.Range("H122").Value = 3264.3794
$ws.Range("I122").Value = 2168.7058
$ws.Range("J122").Value = 4816.5835
$ws.Range("K122").Value = 6506.117400000001
$ws.Range("L122").Value = 14449.7505
$ws.Range("M122").Value = -4056.117400000001
$ws.Range("N122").Value = -19349.7505

$ws.Range("H126").Value = 4782.1377
$ws.Range("I126").Value = 4247.05
$ws.Range("K126").Value = 12741.15
$ws.Range("M126").Value = -10271.15

$ws.Range("H132").Value = 10004109
$ws.Range("I132").Value = 13516142
$ws.Range("J132").Value = 8323.691999999999
$ws.Range("K132").Value = 40548426
$ws.Range("L132").Value = 24971.076
$ws.Range("M132").Value = -40545896
$ws.Range("N132").Value = -30031.076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 36405.547
$ws.Range("J45").Value = 37385.445
$ws.Range("L45").Value = 37385.445
$ws.Range("N45").Value = -38367.445

$ws.Range("H52").Value = 9350
$ws.Range("I52").Value = 9350
$ws.Range("K52").Value = 9350
$ws.Range("M52").Value = -9124

$ws.Range("H62").Value = 43691.625
$ws.Range("I62").Value = 53953.844
$ws.Range("K62").Value = 53953.844
$ws.Range("M62").Value = -53329.844

$ws.Range("H65").Value = 43691.625
$ws.Range("I65").Value = 53953.844
$ws.Range("K65").Value = 269769.22
$ws.Range("M65").Value = -266649.22

$ws.Range("H107").Value = 15873775
$ws.Range("I107").Value = 611.8333
$ws.Range("J107").Value = 37037990
$ws.Range("K107").Value = 1835.4999
$ws.Range("L107").Value = 111113970
$ws.Range("M107").Value = 84.50009999999997
$ws.Range("N107").Value = -111117810

$ws.Range("H122").Value = 172423.92
$ws.Range("I122").Value = 404982.4
$ws.Range("K122").Value = 1214947.2
$ws.Range("M122").Value = -1212497.2

$ws.Range("H132").Value = 3019.5925
$ws.Range("I132").Value = 2825.2932
$ws.Range("J132").Value = 3509.5652
$ws.Range("K132").Value = 8475.8796
$ws.Range("L132").Value = 10528.6956
$ws.Range("M132").Value = -5945.8796
$ws.Range("N132").Value = -15588.6956

$ws.Range("H137").Value = 89999.25
$ws.Range("J137").Value = 89999.25
$ws.Range("L137").Value = 89999.25
$ws.Range("N137").Value = -100199.25

